$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New row 15: store address / city / manager-name / manager-document ---
$ws.Range("A15").Value = "Rua Décio Vilares, 406 - Copacabana, Rio de Janeiro - RJ, 22045-050"
$ws.Range("B15").Value = "Rio de Janeiro"
$ws.Range("C15").Value = "Alon Pinheiro"
$ws.Range("D15").Value = "28456321X"

# --- Formatting for the new row: centered, thin left/right borders only ---
foreach ($col in @("A", "B", "C", "D")) {
    $cell = $ws.Range($col + "15")
    $cell.Borders.Item(7).LineStyle = 1    # xlEdgeLeft / xlContinuous
    $cell.Borders.Item(7).Weight = 2       # xlThin
    $cell.Borders.Item(10).LineStyle = 1   # xlEdgeRight / xlContinuous
    $cell.Borders.Item(10).Weight = 2      # xlThin
    $cell.HorizontalAlignment = -4108      # xlCenter
}

# --- Move the active selection to E15, matching the saved view state ---
$ws.Range("E15").Select()
